$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns A (Id) and B (Title) - shifts Author..Nb_Page left by 2
$ws.Range("A1:B1").EntireColumn.Delete()

# Update the view: clear the split/topLeftCell, move selection to K7
$ws.Range("K7").Select()
